# Updated cryptos list on Fri Jun  7 15:37:48 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "Price" (column D) cell to a literal text value, forcing a
# text number-format first so Excel does not silently reinterpret values
# such as "700.24" or "0.0000250" as numbers (which would drop trailing
# zeros / switch to scientific notation and lose the exact original text).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.003.86"
$ws.Range("E2").Value = "  -0.76%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.799.07"
$ws.Range("E3").Value = "  -1.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "700.24"
$ws.Range("E5").Value = "  -1.21%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "169.76"
$ws.Range("E6").Value = "  -2.04%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.799.25"
$ws.Range("E7").Value = "  -1.11%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.17%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.51%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.12%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +2.13%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.481"
$ws.Range("E12").Value = "  +4.37%  "

# Row 13 - ShibaInu
Set-TextValue $ws.Range("D13") "0.0000250"
$ws.Range("E13").Value = "  -2.79%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "36.22"
$ws.Range("E14").Value = "  -2.33%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.446.83"
$ws.Range("E15").Value = "  -3.91%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.818.69"
$ws.Range("E16").Value = "  -0.73%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "71.144.05"
$ws.Range("E17").Value = "  -0.55%  "

# Row 18 - Chainlink
Set-TextValue $ws.Range("D18") "17.62"
$ws.Range("E18").Value = "  +0.73%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "7.22"
$ws.Range("E19").Value = "  -0.65%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.25%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "512.33"
$ws.Range("E21").Value = "  +2.64%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "10.46"
$ws.Range("E22").Value = "  -2.75%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -2.59%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "83.59"
$ws.Range("E24").Value = "  -2.40%  "

# Row 25 - PEPE
Set-TextValue $ws.Range("D25") "0.0000141"
$ws.Range("E25").Value = "  -3.74%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D26") "12.67"
$ws.Range("E26").Value = "  +3.48%  "

# Row 27 - WrappedeETH
$ws.Range("D27").Value = "3.953.12"
$ws.Range("E27").Value = "  -1.10%  "

# Row 28 - RenderToken
Set-TextValue $ws.Range("D28") "10.26"
$ws.Range("E28").Value = "  -4.59%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.07%  "

# Row 30 - Fetch.AI
$ws.Range("E30").Value = "  -5.65%  "

# Row 31 - PancakeSwap
Set-TextValue $ws.Range("D31") "2.99"
$ws.Range("E31").Value = "  -4.73%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +1.05%  "

# Row 33 - NEARProtocol
Set-TextValue $ws.Range("D33") "7.29"
$ws.Range("E33").Value = "  -3.20%  "

# Row 34 - EthereumClassic
Set-TextValue $ws.Range("D34") "29.10"
$ws.Range("E34").Value = "  -1.35%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  -4.82%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  +0.58%  "

# Row 37 - RenzoRestakedETH
$ws.Range("D37").Value = "3.769.32"
$ws.Range("E37").Value = "  -1.02%  "

# Row 38 - Binance-PegBSC-USD
$ws.Range("E38").Value = "  +0.25%  "

# Row 39 - Filecoin
Set-TextValue $ws.Range("D39") "6.64"
$ws.Range("E39").Value = "  +9.93%  "

# Row 40 - Hedera
$ws.Range("E40").Value = "  -2.71%  "

# Row 41 - Stacks
Set-TextValue $ws.Range("D41") "2.35"
$ws.Range("E41").Value = "  +0.61%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  -2.64%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  +0.00%  "

# Row 44 - dogwifhat
Set-TextValue $ws.Range("D44") "3.18"
$ws.Range("E44").Value = "  -6.26%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.06%  "

# Row 46 - Monero
Set-TextValue $ws.Range("D46") "164.31"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47 - was Bittensor, now OKB
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D47") "49.34"
$ws.Range("E47").Value = "  +0.51%  "

# Row 48 - was OKB, now FLOKI
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D48") "0.000305"
$ws.Range("E48").Value = "  -4.61%  "

# Row 49 - was FLOKI, now Bittensor
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D49") "429.79"
$ws.Range("E49").Value = "  -0.42%  "

# Row 50 - Cosmos (unchanged name/link)
Set-TextValue $ws.Range("D50") "8.63"
$ws.Range("E50").Value = "  -1.39%  "

# Row 51 - was InjectiveProtocol, now TheGraph
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D51") "0.296"
$ws.Range("E51").Value = "  -1.35%  "
